# Adds season-record columns (Wins / Losses / Ties) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new column headers in AD1:AF1 ---
# Copy the formatting of the existing header cell (AC1, style s="1") onto the
# new header cells before writing their text, so they match the rest of the
# header row exactly.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (rows 2-47): season record values for every player row ---
$wins = 66
$losses = 96
$ties = 0

for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # column AD
    $ws.Cells.Item($r, 31).Value = $losses  # column AE
    $ws.Cells.Item($r, 32).Value = $ties    # column AF
}
